$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.248835333333334
$ws.Range("H2").Value = 6.746506
$ws.Range("I2").Value = 0.03590294220158827
$ws.Range("J2").Value = 0.03590294220158827
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.737936333333333
$ws.Range("N2").Value = 5.213808999999999
$ws.Range("O2").Value = 0.01383107950678261
$ws.Range("P2").Value = 0.01383107950678261
$ws.Range("Q2").Value = 3.908332633483778
$ws.Range("R2").Value = 35.174993701354
$ws.Range("S2").Value = 0.0004965764481175882
$ws.Range("T2").Value = 0.0004965764481175882
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.248835333333334
$ws.Range("H3").Value = 6.746506
$ws.Range("I3").Value = 0.03590294220158827
$ws.Range("J3").Value = 0.03590294220158827
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.524875
$ws.Range("N3").Value = 10.574625
$ws.Range("O3").Value = 0.02805213599681367
$ws.Range("P3").Value = 0.02805213599681367
$ws.Range("Q3").Value = 7.926863445583334
$ws.Range("R3").Value = 71.34177101025
$ws.Range("S3").Value = 0.001007154217324695
$ws.Range("T3").Value = 0.001007154217324695
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.248835333333334
$ws.Range("H4").Value = 6.746506
$ws.Range("I4").Value = 0.03590294220158827
$ws.Range("J4").Value = 0.03590294220158827
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 49.07777633333333
$ws.Range("N4").Value = 147.233329
$ws.Range("O4").Value = 0.3905773839140027
$ws.Range("P4").Value = 0.3905773839140027
$ws.Range("Q4").Value = 110.3678374998305
$ws.Range("R4").Value = 993.310537498474
$ws.Range("S4").Value = 0.01402287723991199
$ws.Range("T4").Value = 0.01402287723991199
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.248835333333334
$ws.Range("H5").Value = 6.746506
$ws.Range("I5").Value = 0.03590294220158827
$ws.Range("J5").Value = 0.03590294220158827
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 69.40412133333334
$ws.Range("N5").Value = 208.212364
$ws.Range("O5").Value = 0.5523412462518597
$ws.Range("P5").Value = 0.5523412462518597
$ws.Range("Q5").Value = 156.0784403333538
$ws.Range("R5").Value = 1404.705963000184
$ws.Range("S5").Value = 0.01983067583973375
$ws.Range("T5").Value = 0.01983067583973375
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.248835333333334
$ws.Range("H6").Value = 6.746506
$ws.Range("I6").Value = 0.03590294220158827
$ws.Range("J6").Value = 0.03590294220158827
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.909715333333333
$ws.Range("N6").Value = 5.729146
$ws.Range("O6").Value = 0.01519815433054137
$ws.Range("P6").Value = 0.01519815433054137
$ws.Range("Q6").Value = 4.294635318208445
$ws.Range("R6").Value = 38.651717863876
$ws.Range("S6").Value = 0.0005456584565002455
$ws.Range("T6").Value = 0.0005456584565002455
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 44.29005966666667
$ws.Range("H7").Value = 132.870179
$ws.Range("I7").Value = 0.7070964373190639
$ws.Range("J7").Value = 0.7070964373190639
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.737936333333333
$ws.Range("N7").Value = 5.213808999999999
$ws.Range("O7").Value = 0.01383107950678261
$ws.Range("P7").Value = 0.01383107950678261
$ws.Range("Q7").Value = 76.97330390020122
$ws.Range("R7").Value = 692.759735101811
$ws.Range("S7").Value = 0.0097799070435227
$ws.Range("T7").Value = 0.0097799070435227
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 44.29005966666667
$ws.Range("H8").Value = 132.870179
$ws.Range("I8").Value = 0.7070964373190639
$ws.Range("J8").Value = 0.7070964373190639
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.524875
$ws.Range("N8").Value = 10.574625
$ws.Range("O8").Value = 0.02805213599681367
$ws.Range("P8").Value = 0.02805213599681367
$ws.Range("Q8").Value = 156.1169240675417
$ws.Range("R8").Value = 1405.052316607875
$ws.Range("S8").Value = 0.01983556542253681
$ws.Range("T8").Value = 0.01983556542253681
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 44.29005966666667
$ws.Range("H9").Value = 132.870179
$ws.Range("I9").Value = 0.7070964373190639
$ws.Range("J9").Value = 0.7070964373190639
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 49.07777633333333
$ws.Range("N9").Value = 147.233329
$ws.Range("O9").Value = 0.3905773839140027
$ws.Range("P9").Value = 0.3905773839140027
$ws.Range("Q9").Value = 2173.657642110655
$ws.Range("R9").Value = 19562.91877899589
$ws.Range("S9").Value = 0.2761758766629915
$ws.Range("T9").Value = 0.2761758766629915
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 44.29005966666667
$ws.Range("H10").Value = 132.870179
$ws.Range("I10").Value = 0.7070964373190639
$ws.Range("J10").Value = 0.7070964373190639
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 69.40412133333334
$ws.Range("N10").Value = 208.212364
$ws.Range("O10").Value = 0.5523412462518597
$ws.Range("P10").Value = 0.5523412462518597
$ws.Range("Q10").Value = 3073.912674965907
$ws.Range("R10").Value = 27665.21407469316
$ws.Range("S10").Value = 0.3905585274090617
$ws.Range("T10").Value = 0.3905585274090617
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 44.29005966666667
$ws.Range("H11").Value = 132.870179
$ws.Range("I11").Value = 0.7070964373190639
$ws.Range("J11").Value = 0.7070964373190639
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.909715333333333
$ws.Range("N11").Value = 5.729146
$ws.Range("O11").Value = 0.01519815433054137
$ws.Range("P11").Value = 0.01519815433054137
$ws.Range("Q11").Value = 84.58140605968157
$ws.Range("R11").Value = 761.2326545371341
$ws.Range("S11").Value = 0.01074656078095111
$ws.Range("T11").Value = 0.01074656078095111
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 16.09762433333333
$ws.Range("H12").Value = 48.292873
$ws.Range("I12").Value = 0.2570006204793478
$ws.Range("J12").Value = 0.2570006204793479
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 1.737936333333333
$ws.Range("N12").Value = 5.213808999999999
$ws.Range("O12").Value = 0.01383107950678261
$ws.Range("P12").Value = 0.01383107950678261
$ws.Range("Q12").Value = 27.97664620925077
$ws.Range("R12").Value = 251.789815883257
$ws.Range("S12").Value = 0.003554596015142323
$ws.Range("T12").Value = 0.003554596015142324
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 16.09762433333333
$ws.Range("H13").Value = 48.292873
$ws.Range("I13").Value = 0.2570006204793478
$ws.Range("J13").Value = 0.2570006204793479
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.524875
$ws.Range("N13").Value = 10.574625
$ws.Range("O13").Value = 0.02805213599681367
$ws.Range("P13").Value = 0.02805213599681367
$ws.Range("Q13").Value = 56.74211357195833
$ws.Range("R13").Value = 510.6790221476249
$ws.Range("S13").Value = 0.007209416356952161
$ws.Range("T13").Value = 0.007209416356952163
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 16.09762433333333
$ws.Range("H14").Value = 48.292873
$ws.Range("I14").Value = 0.2570006204793478
$ws.Range("J14").Value = 0.2570006204793479
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 49.07777633333333
$ws.Range("N14").Value = 147.233329
$ws.Range("O14").Value = 0.3905773839140027
$ws.Range("P14").Value = 0.3905773839140027
$ws.Range("Q14").Value = 790.0356065293573
$ws.Range("R14").Value = 7110.320458764217
$ws.Range("S14").Value = 0.1003786300110991
$ws.Range("T14").Value = 0.1003786300110991
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 16.09762433333333
$ws.Range("H15").Value = 48.292873
$ws.Range("I15").Value = 0.2570006204793478
$ws.Range("J15").Value = 0.2570006204793479
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 69.40412133333334
$ws.Range("N15").Value = 208.212364
$ws.Range("O15").Value = 0.5523412462518597
$ws.Range("P15").Value = 0.5523412462518597
$ws.Range("Q15").Value = 1117.241472409086
$ws.Range("R15").Value = 10055.17325168177
$ws.Range("S15").Value = 0.1419520430030642
$ws.Range("T15").Value = 0.1419520430030642
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 16.09762433333333
$ws.Range("H16").Value = 48.292873
$ws.Range("I16").Value = 0.2570006204793478
$ws.Range("J16").Value = 0.2570006204793479
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 1.909715333333333
$ws.Range("N16").Value = 5.729146
$ws.Range("O16").Value = 0.01519815433054137
$ws.Range("P16").Value = 0.01519815433054137
$ws.Range("Q16").Value = 30.74188001960644
$ws.Range("R16").Value = 276.676920176458
$ws.Range("S16").Value = 0.00390593509309002
$ws.Range("T16").Value = 0.00390593509309002
